# Swap the differing fields between row 2 and row 3 (the two occurrence
# records got their row order flipped; most columns are identical between
# the two rows, so only the columns that actually differ need touching).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$columns = @("A", "B", "E", "F", "G", "H", "P", "Q", "R", "S", "AC")

foreach ($col in $columns) {
    $cell2 = $ws.Range($col + "2")
    $cell3 = $ws.Range($col + "3")

    $v2 = $cell2.Value2
    $v3 = $cell3.Value2

    $cell2.Value2 = $v3
    $cell3.Value2 = $v2
}
